$p = $ppt.ActivePresentation

# ppPlaceholderDate = 16
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.Type -eq 14) {
            $pf = $shp.PlaceholderFormat
            if ($pf.Type -eq $ppPlaceholderDate) {
                $shp.TextFrame.TextRange.Text = "1/10/2021"
            }
        }
    }
}

# Slide master's own Date placeholder
$master = $p.Slides.Item(1).Master
Update-DatePlaceholder $master.Shapes

# Every slide (custom) layout's Date placeholder
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}
